$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Row 22 ("Kaskaden auf dem alten Friedhof"): the "Interest additional" (N) field becomes a
# multi-choice combination of "Museum" + "City Experience", and the From/To hours (K/L) get set.
$ws.Range("N22").Value = "Museum, City Experience"
$ws.Range("K22").Value = 10
$ws.Range("L22").Value = 15

# Column N widens to fit the longer selected text.
$ws.Columns.Item(14).ColumnWidth = 47

# The active selection moves to N26 (below the edited cell).
[void]$ws.Range("N26").Select()
